$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reading data from a specific cell
$name = $ws.Range("A1").Value2

# Inserting data in specific cells (column B, rows 1-4)
$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 1
